# Correlation Column Logic and Data Gen.
# Adds a GEN_ORDER column (between ROWS and COLUMN_ORDER), reorders/extends
# the table rows, and appends a new S_SUPPL_PARTNER row plus a trailing
# blank row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new column D ("GEN_ORDER"). This shifts the old D column
#    (COLUMN_ORDER) to E, carrying its values/styles along with it.
# ---------------------------------------------------------------------
$ws.Columns("D").Insert()

# New column D should look like the numeric "ROWS" column (right aligned,
# #,##0 number format, no visible border) - copy that formatting down the
# whole column (header included, fixed up right after).
$ws.Range("C2").Copy()
$ws.Range("D1:D8").PasteSpecial(-4122)

# Header row ------------------------------------------------------------
$ws.Range("D1").Value = "GEN_ORDER"
$ws.Range("E1").Value = "COLUMN_ORDER"

# ---------------------------------------------------------------------
# 2. Grow the table by two rows (row 9 for the new S_SUPPL_PARTNER entry,
#    row 10 as a new trailing blank row) by copying row 8's formatting
#    down before any values are written, so the new rows inherit the
#    correct per-cell styles (and row height) instead of column defaults.
# ---------------------------------------------------------------------
$ws.Range("A8:E8").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)
$ws.Rows(9).RowHeight = 19.5

$ws.Range("A9:E9").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)
$ws.Rows(10).RowHeight = 19.5

# Column E (old D, "COLUMN_ORDER" text) keeps its own distinct text style
# (right-aligned, General number format) the whole way down - row 8's
# copy-down above carries the numeric ROWS/GEN_ORDER look instead, so
# reapply E2's format onto E3:E10.
$ws.Range("E2").Copy()
$ws.Range("E3:E10").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Re-key the data rows. Row 2 (S_SUPPL_GEN) keeps its place; the rest
#    are reordered, a new S_SUPPL_PARTNER row is appended, and a new
#    trailing blank row is added - matching the "GEN_ORDER" sequence.
# ---------------------------------------------------------------------
$domain = "vendor"

$rows = @(
    @{ table = "S_SUPPL_GEN";      rows = 15; order = 1; cols = "LIFNR, BU_GROUP, KTOKK, NAME_FIRST, NAME_FIRST_P, NAME_LAST_P, BPEXT, STREET, POST_CODE1, CITY1, COUNTRY, REGION, LANGU_CORR, TELNR_LONG, SMTP_ADDR" },
    @{ table = "S_LFA1_TEXT";      rows = 5;  order = 2; cols = "LIFNR,TDSPRAS" },
    @{ table = "S_SUPPL_ADDR";     rows = 5;  order = 3; cols = "LIFNR" },
    @{ table = "S_LFM1_TEXT";      rows = 5;  order = 4; cols = "LIFNR" },
    @{ table = "S_SUPPL_WITH_TAX"; rows = 5;  order = 5; cols = "LIFNR,BUKRS" },
    @{ table = "S_ROLES";          rows = 5;  order = 6; cols = "LIFNR,BP_ROLE" },
    @{ table = "S_ADDR_USAGE";     rows = 5;  order = 7; cols = "LIFNR, ADR_KIND" },
    @{ table = "S_SUPPL_PARTNER";  rows = 5;  order = 8; cols = "LIFNR,EKORG,PARVW,LIFN2,DEFPA" }
)

$r = 2
foreach ($rec in $rows) {
    $ws.Cells.Item($r, 1).Value = $domain
    $ws.Cells.Item($r, 2).Value = $rec.table
    $ws.Cells.Item($r, 3).Value = $rec.rows
    $ws.Cells.Item($r, 4).Value = $rec.order
    $ws.Cells.Item($r, 5).Value = $rec.cols
    $r = $r + 1
}

# Row 10 stays fully blank (trailing spacer row), matching row 9's formats.
$ws.Range("A10").Value = ""
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""

# ---------------------------------------------------------------------
# 4. Column widths - tighten B/C/D/E to the new layout.
# ---------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 16.72
$ws.Columns("C").ColumnWidth = 108.58
$ws.Columns("D").ColumnWidth = 12.43
$ws.Columns("E").ColumnWidth = 141.57

$ws.Range("A1").Select() | Out-Null
